$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.539.75"
$ws.Range("E2").Value = "  +1.33%  "
$ws.Range("D3").Value = "1.882.86"
$ws.Range("E3").Value = "  +1.34%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "248.34"
$ws.Range("E5").Value = "  +6.67%  "
$ws.Range("E6").Value = "  -0.05%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4760"
$ws.Range("E7").Value = "  +1.35%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2933"
$ws.Range("E8").Value = "  +3.69%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06531"
$ws.Range("E9").Value = "  +1.33%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "22.04"
$ws.Range("E10").Value = "  +6.29%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "98.14"
$ws.Range("E11").Value = "  +5.26%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07733"
$ws.Range("E12").Value = "  +0.91%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.7409"
$ws.Range("E13").Value = "  +9.74%  "
$ws.Range("D14").Value = "1.885.24"
$ws.Range("E14").Value = "  +1.58%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.151"
$ws.Range("E15").Value = "  +2.40%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "274.22"
$ws.Range("E16").Value = "  +2.81%  "
$ws.Range("D17").Value = "30.543.79"
$ws.Range("E17").Value = "  +1.47%  "
$ws.Range("E18").Value = "  +1.85%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007576"
$ws.Range("E19").Value = "  +1.17%  "
$ws.Range("D21").Value = "2.135.39"
$ws.Range("E21").Value = "  +2.10%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.001"
$ws.Range("E22").Value = "  +0.00%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.260"
$ws.Range("E23").Value = "  +2.58%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.205"
$ws.Range("E24").Value = "  +2.37%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.297"
$ws.Range("E25").Value = "  +0.64%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "163.77"
$ws.Range("E26").Value = "  -0.87%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.92"
$ws.Range("E27").Value = "  +2.47%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.946"
$ws.Range("E28").Value = "  +3.94%  "
$ws.Range("E29").Value = "  +2.89%  "
$ws.Range("E30").Value = "  +0.06%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.522"
$ws.Range("E31").Value = "  +5.29%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.334"
$ws.Range("E32").Value = "  +3.89%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.113"
$ws.Range("E33").Value = "  +3.94%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04826"
$ws.Range("E34").Value = "  +4.25%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.132"
$ws.Range("E35").Value = "  +2.32%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7022"
$ws.Range("E36").Value = "  +3.21%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.713"
$ws.Range("E37").Value = "  -0.01%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01869"
$ws.Range("E38").Value = "  +3.40%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.753"
$ws.Range("E39").Value = "  +1.64%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.320"
$ws.Range("E40").Value = "  +1.08%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.000"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "71.62"
$ws.Range("E42").Value = "  +2.34%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.4224"
$ws.Range("E43").Value = "  +4.96%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.8449"
$ws.Range("E44").Value = "  +2.10%  "
$ws.Range("E45").Value = "  -0.04%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "102.94"
$ws.Range("E46").Value = "  +0.82%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.307"
$ws.Range("E47").Value = "  +2.18%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.105"
$ws.Range("E48").Value = "  +3.49%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "35.65"
$ws.Range("E49").Value = "  +4.78%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "918.20"
$ws.Range("E50").Value = "  -0.10%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.3902"
$ws.Range("E51").Value = "  +4.78%  "
